# Reorder the Data sheet rows (user-data filename fix).
# Only the row order changes; cell contents per id are unchanged, so we
# only write a cell when the new value differs from what currently sits
# in that grid position.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 2: -> DZ
$ws.Cells.Item(2, 1).Value = 'DZ'
$ws.Cells.Item(2, 2).Value = 'Algeria'
$ws.Cells.Item(2, 3).Value = 'Sahrawi refugee crisis'
$ws.Cells.Item(2, 5).Value = 13
$ws.Cells.Item(2, 6).Value = '9 or over'

# Row 3: -> BD
$ws.Cells.Item(3, 1).Value = 'BD'
$ws.Cells.Item(3, 2).Value = 'Bangladesh'
$ws.Cells.Item(3, 3).Value = 'Chittagong Hill Tracts'
$ws.Cells.Item(3, 5).Value = 9
$ws.Cells.Item(3, 6).Value = '6 or over'

# Row 4: -> BF
$ws.Cells.Item(4, 1).Value = 'BF'
$ws.Cells.Item(4, 2).Value = 'Burkina Faso'
$ws.Cells.Item(4, 3).Value = 'Sahel regional crisis'
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = '1 or over'

# Row 5: -> CM
$ws.Cells.Item(5, 1).Value = 'CM'
$ws.Cells.Item(5, 2).Value = 'Cameroon'
$ws.Cells.Item(5, 3).Value = ""
$ws.Cells.Item(5, 5).Value = 2
$ws.Cells.Item(5, 6).Value = '1 or over'

# Row 6: -> CF
$ws.Cells.Item(6, 1).Value = 'CF'
$ws.Cells.Item(6, 2).Value = 'CAR'
$ws.Cells.Item(6, 3).Value = 'Internal armed conflict'
$ws.Cells.Item(6, 5).Value = 5
$ws.Cells.Item(6, 6).Value = '3 or over'

# Row 7: -> TD
$ws.Cells.Item(7, 1).Value = 'TD'
$ws.Cells.Item(7, 2).Value = 'Chad'
$ws.Cells.Item(7, 3).Value = 'Sahel regional crisis'
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = '3 or over'

# Row 8: -> CO
$ws.Cells.Item(8, 1).Value = 'CO'
$ws.Cells.Item(8, 2).Value = 'Colombia'
$ws.Cells.Item(8, 3).Value = 'Armed conflict'
$ws.Cells.Item(8, 5).Value = 10

# Row 9: -> CD
$ws.Cells.Item(9, 1).Value = 'CD'
$ws.Cells.Item(9, 2).Value = 'DRC'
$ws.Cells.Item(9, 3).Value = 'Affected by humanitarian crisis caused by LRA'
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = '3 or over'

# Row 10: -> EC
$ws.Cells.Item(10, 1).Value = 'EC'
$ws.Cells.Item(10, 2).Value = 'Ecuador'
$ws.Cells.Item(10, 3).Value = 'Colombian refugees'

# Row 11: -> EG
$ws.Cells.Item(11, 1).Value = 'EG'
$ws.Cells.Item(11, 2).Value = 'Egypt'
$ws.Cells.Item(11, 3).Value = 'Refugee crisis'

# Row 12: -> GE
$ws.Cells.Item(12, 1).Value = 'GE'
$ws.Cells.Item(12, 2).Value = 'Georgia'
$ws.Cells.Item(12, 3).Value = 'Abkhazia'

# Row 13: -> GN
$ws.Cells.Item(13, 1).Value = 'GN'
$ws.Cells.Item(13, 2).Value = 'Guinea'
$ws.Cells.Item(13, 5).Value = 1

# Row 14: -> HT
$ws.Cells.Item(14, 1).Value = 'HT'
$ws.Cells.Item(14, 2).Value = 'Haiti'
$ws.Cells.Item(14, 5).Value = 2
$ws.Cells.Item(14, 6).Value = '1 or over'

# Row 15: -> IN

# Row 16: -> ID
$ws.Cells.Item(16, 1).Value = 'ID'
$ws.Cells.Item(16, 2).Value = 'Indonesia'
$ws.Cells.Item(16, 3).Value = ""
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = '3 or over'

# Row 17: -> KE
$ws.Cells.Item(17, 1).Value = 'KE'
$ws.Cells.Item(17, 2).Value = 'Kenya'
$ws.Cells.Item(17, 3).Value = 'Somali refugee crisis'
$ws.Cells.Item(17, 5).Value = 1

# Row 18: -> LY

# Row 19: -> ML

# Row 20: -> MR
$ws.Cells.Item(20, 1).Value = 'MR'
$ws.Cells.Item(20, 2).Value = 'Mauritania'
$ws.Cells.Item(20, 3).Value = 'Sahel regional crisis'
$ws.Cells.Item(20, 5).Value = 1
$ws.Cells.Item(20, 6).Value = '1 or over'

# Row 21: -> MM
$ws.Cells.Item(21, 1).Value = 'MM'
$ws.Cells.Item(21, 2).Value = 'Myanmar'
$ws.Cells.Item(21, 3).Value = 'Northern Rakhine State and Kachin and Shan State conflict'
$ws.Cells.Item(21, 5).Value = 13
$ws.Cells.Item(21, 6).Value = '9 or over'

# Row 22: -> NP
$ws.Cells.Item(22, 1).Value = 'NP'
$ws.Cells.Item(22, 2).Value = 'Nepal'
$ws.Cells.Item(22, 3).Value = 'Bhutanese refugees'
$ws.Cells.Item(22, 5).Value = 8
$ws.Cells.Item(22, 6).Value = '6 or over'

# Row 23: -> NE
$ws.Cells.Item(23, 1).Value = 'NE'
$ws.Cells.Item(23, 2).Value = 'Niger'
$ws.Cells.Item(23, 3).Value = 'Sahel regional crisis'
$ws.Cells.Item(23, 5).Value = 1
$ws.Cells.Item(23, 6).Value = '1 or over'

# Row 24: -> PK
$ws.Cells.Item(24, 1).Value = 'PK'
$ws.Cells.Item(24, 2).Value = 'Pakistan'
$ws.Cells.Item(24, 3).Value = 'Conflict, IDP crisis'
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = '3 or over'

# Row 25: -> PG
$ws.Cells.Item(25, 1).Value = 'PG'
$ws.Cells.Item(25, 2).Value = 'Papua New Guinea'
$ws.Cells.Item(25, 3).Value = ""
$ws.Cells.Item(25, 5).Value = 1

# Row 26: -> PH
$ws.Cells.Item(26, 1).Value = 'PH'
$ws.Cells.Item(26, 2).Value = 'Philippines'
$ws.Cells.Item(26, 3).Value = 'Mindanao crisis'
$ws.Cells.Item(26, 5).Value = 2
$ws.Cells.Item(26, 6).Value = '1 or over'

# Row 27: -> RU

# Row 28: -> SO
$ws.Cells.Item(28, 1).Value = 'SO'
$ws.Cells.Item(28, 2).Value = 'Somalia'
$ws.Cells.Item(28, 3).Value = ""
$ws.Cells.Item(28, 5).Value = 2
$ws.Cells.Item(28, 6).Value = '1 or over'

# Row 29: -> LK
$ws.Cells.Item(29, 1).Value = 'LK'
$ws.Cells.Item(29, 2).Value = 'Sri Lanka'
$ws.Cells.Item(29, 3).Value = 'Returning IDPs'

# Row 30: -> SD
$ws.Cells.Item(30, 1).Value = 'SD'
$ws.Cells.Item(30, 2).Value = 'Sudan'
$ws.Cells.Item(30, 3).Value = 'Darfur, refugees, transitional areas, East Sudan'
$ws.Cells.Item(30, 5).Value = 5

# Row 31: -> TJ
$ws.Cells.Item(31, 1).Value = 'TJ'
$ws.Cells.Item(31, 2).Value = 'Tajikistan'
$ws.Cells.Item(31, 3).Value = ""
$ws.Cells.Item(31, 5).Value = 1
$ws.Cells.Item(31, 6).Value = '1 or over'

# Row 32: -> TZ
$ws.Cells.Item(32, 1).Value = 'TZ'
$ws.Cells.Item(32, 2).Value = 'Tanzania'

# Row 33: -> TH
$ws.Cells.Item(33, 1).Value = 'TH'
$ws.Cells.Item(33, 2).Value = 'Thailand'
$ws.Cells.Item(33, 3).Value = 'Burmese border'
$ws.Cells.Item(33, 5).Value = 7
$ws.Cells.Item(33, 6).Value = '6 or over'

# Row 34: -> UG

# Row 35: -> VE

# Row 36: -> YE
